{"js": "// Update the \"Historial de Versiones\" table:\n//  1. Change the date of version \"01\" from 13/02/2018 to 25/03/2019.\n//  2. Append a new row documenting version \"02\" (02/04/2019 - Ejecuci\u00f3n de\n//     casos de prueba).\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst versionTable = tables.items[0];\n\n// 1) Update the date cell text (01 -> 25/03/2019) while keeping formatting.\nconst dateResults = versionTable.search(\"13/02/2018\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\n\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"25/03/2019\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) Add the new version row at the end of the table.\nversionTable.addRows(Word.InsertLocation.end, 1, [\n  [\"02\", \"02/04/2019\", \"Ejecuci\u00f3n de casos de prueba\", \"\"]\n]);\nawait context.sync();\n", "ps1": "# Update the \"Historial de Versiones\" table:\n#  1. Change the date of version \"01\" from 13/02/2018 to 25/03/2019.\n#  2. Append a new row documenting version \"02\" (02/04/2019 - Ejecuci\u00f3n de\n#     casos de prueba).\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# 1) Update the date cell text (01 -> 25/03/2019) while keeping formatting.\n$find = $t.Range.Find\n$find.Execute(\"13/02/2018\", $false, $false, $false, $false, $false, $true, 1, $false, \"25/03/2019\", 2) | Out-Null\n\n# 2) Add the new version row at the end of the table.\n$newRow = $t.Rows.Add()\n$newRow.Cells.Item(1).Range.Text = \"02\"\n$newRow.Cells.Item(2).Range.Text = \"02/04/2019\"\n$newRow.Cells.Item(3).Range.Text = \"Ejecuci\u00f3n de casos de prueba\"\n"}
